$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (old first data row) - this shifts all subsequent rows
# (3..33) up by one, so that the Month (A) and Actual (B) values line up
# with the new row numbers, and removes the old trailing row 33.
$ws.Rows(2).Delete()

# Update the Predicted (column C) values for all 32 remaining data rows
# with the newly retrained model's predictions.
$ws.Range("C2").Value = 5.263088226318359
$ws.Range("C3").Value = 5.263363838195801
$ws.Range("C4").Value = 5.263692855834961
$ws.Range("C5").Value = 5.264224529266357
$ws.Range("C6").Value = 5.26333475112915
$ws.Range("C7").Value = 5.262760639190674
$ws.Range("C8").Value = 5.263341903686523
$ws.Range("C9").Value = 5.262135982513428
$ws.Range("C10").Value = 5.26285982131958
$ws.Range("C11").Value = 5.262763500213623
$ws.Range("C12").Value = 5.26339054107666
$ws.Range("C13").Value = 5.26413106918335
$ws.Range("C14").Value = 5.264116764068604
$ws.Range("C15").Value = 5.264317035675049
$ws.Range("C16").Value = 5.264308452606201
$ws.Range("C17").Value = 5.263980388641357
$ws.Range("C18").Value = 5.26328706741333
$ws.Range("C19").Value = 5.262651443481445
$ws.Range("C20").Value = 5.26248025894165
$ws.Range("C21").Value = 5.260635852813721
$ws.Range("C22").Value = 5.261858463287354
$ws.Range("C23").Value = 5.262577056884766
$ws.Range("C24").Value = 5.263165950775146
$ws.Range("C25").Value = 5.263201236724854
$ws.Range("C26").Value = 5.263431072235107
$ws.Range("C27").Value = 5.263772964477539
$ws.Range("C28").Value = 5.263920307159424
$ws.Range("C29").Value = 5.263461589813232
$ws.Range("C30").Value = 5.262987613677979
$ws.Range("C31").Value = 5.262224674224854
$ws.Range("C32").Value = 5.262828350067139
